# Apply cryptos list price/volume updates per commit
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.751.73"
$ws.Range("E2").Value = "  -0.07%  "

$ws.Range("D3").Value = "2.731.73"
$ws.Range("E3").Value = "  -0.48%  "

$ws.Range("E4").Value = "  -0.04%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "561.59"
$ws.Range("E5").Value = "  -1.99%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "160.12"
$ws.Range("E6").Value = "  +1.96%  "

$ws.Range("E7").Value = "  +0.01%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.594"
$ws.Range("E8").Value = "  -1.25%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.108"
$ws.Range("E9").Value = "  -0.49%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.167"
$ws.Range("E10").Value = "  +3.94%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.55"
$ws.Range("E11").Value = "  -1.58%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.374"
$ws.Range("E12").Value = "  -1.68%  "

$ws.Range("D13").Value = "3.214.45"
$ws.Range("E13").Value = "  -0.64%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "26.72"
$ws.Range("E14").Value = "  +0.64%  "

$ws.Range("D15").Value = "63.599.17"
$ws.Range("E15").Value = "  +0.25%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0000148"
$ws.Range("E16").Value = "  -1.15%  "

$ws.Range("D17").Value = "2.738.61"
$ws.Range("E17").Value = "  -0.45%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "12.28"
$ws.Range("E18").Value = "  +1.48%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.69"
$ws.Range("E19").Value = "  -2.25%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "354.10"
$ws.Range("E20").Value = "  -0.25%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.59"
$ws.Range("E21").Value = "  -1.43%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.999"
$ws.Range("E22").Value = "  +0.13%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.517"
$ws.Range("E23").Value = "  -2.74%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "64.02"
$ws.Range("E24").Value = "  -1.72%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.169"
$ws.Range("E25").Value = "  -0.33%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.00"
$ws.Range("E26").Value = "  +0.18%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.28"
$ws.Range("E27").Value = "  -2.28%  "

$ws.Range("D28").Value = "0.0₃0900"
$ws.Range("E28").Value = "  -0.47%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.97"
$ws.Range("E29").Value = "  +1.49%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.36"
$ws.Range("E30").Value = "  +9.99%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.14"
$ws.Range("E31").Value = "  +1.16%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "167.16"
$ws.Range("E32").Value = "  -1.36%  "

$ws.Range("B33").Value = "ImmutableX"
$ws.Range("C33").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.50"
$ws.Range("E33").Value = "  +3.80%  "

$ws.Range("B34").Value = "NEARProtocol"
$ws.Range("C34").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.90"
$ws.Range("E34").Value = "  +0.31%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "19.98"
$ws.Range("E35").Value = "  -0.82%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.999"
$ws.Range("E36").Value = "  -0.01%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.80"
$ws.Range("E37").Value = "  +0.92%  "

$ws.Range("B38").Value = "Bittensor"
$ws.Range("C38").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "343.73"
$ws.Range("E38").Value = "  +3.73%  "

$ws.Range("B39").Value = "SuiNetwork"
$ws.Range("C39").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.962"
$ws.Range("E39").Value = "  -1.56%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "6.28"
$ws.Range("E40").Value = "  +2.05%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "4.07"
$ws.Range("E41").Value = "  -1.92%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "38.51"
$ws.Range("E42").Value = "  -1.13%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "21.59"
$ws.Range("E43").Value = "  +0.55%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "20.85"
$ws.Range("E44").Value = "  -2.39%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0580"
$ws.Range("E45").Value = "  -0.84%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.630"
$ws.Range("E46").Value = "  +0.98%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0249"
$ws.Range("E47").Value = "  -1.79%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0991"
$ws.Range("E48").Value = "  -1.88%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "132.00"
$ws.Range("E49").Value = "  -2.55%  "

$ws.Range("E50").Value = "  -0.11%  "

$ws.Range("E51").Value = "  +0.33%  "

